# Insert a new weekly price record at the top of the "Poroto verde" data
# block (sheet data starts at row 2; this workbook currently has data
# rows through row 141, dimension A1:R141). The new record is inserted
# as the new row 95, shifting the previous rows 95-141 down to 96-142.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 95:141 down by one row, creating a blank row 95.
$ws.Rows(95).Insert()

# Populate the new row 95 with the latest weekly record.
$ws.Range("A95").Value = 7
$ws.Range("B95").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C95").Value = "Ñuble"
$ws.Range("D95").Value = 45089
$ws.Range("E95").Value = 16
$ws.Range("F95").Value = 100112031
$ws.Range("G95").Value = "Poroto verde"
$ws.Range("H95").Value = "Magnum"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 30
$ws.Range("K95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("M95").Value = 25000
$ws.Range("N95").Value = "`$/malla 25 kilos"
$ws.Range("O95").Value = "Perú"
$ws.Range("P95").Value = 1000
$ws.Range("Q95").Value = 25
$ws.Range("R95").Value = "Hortaliza"
